# Updates cryptos list: refreshed prices/volumes and re-ranked rows 21-25 & 40-41
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.738.75'
$ws.Range('E2').Value = '  +3.49%  '

$ws.Range('D3').Value = '2.448.03'
$ws.Range('E3').Value = '  +2.13%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '''577.93'
$ws.Range('E5').Value = '  +2.96%  '

$ws.Range('D6').Value = '''145.88'
$ws.Range('E6').Value = '  +3.34%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  +0.54%  '

$ws.Range('D9').Value = '2.446.63'
$ws.Range('E9').Value = '  +1.85%  '

$ws.Range('E10').Value = '  +2.42%  '

$ws.Range('E11').Value = '  +1.14%  '

$ws.Range('D12').Value = '''5.25'
$ws.Range('E12').Value = '  +1.72%  '

$ws.Range('D13').Value = '''0.354'
$ws.Range('E13').Value = '  +3.24%  '

$ws.Range('D14').Value = '''28.47'
$ws.Range('E14').Value = '  +9.50%  '

$ws.Range('D15').Value = '''0.0000179'
$ws.Range('E15').Value = '  +6.40%  '

$ws.Range('D16').Value = '2.891.56'
$ws.Range('E16').Value = '  +2.17%  '

$ws.Range('D17').Value = '62.595.56'
$ws.Range('E17').Value = '  +3.54%  '

$ws.Range('D18').Value = '2.454.01'
$ws.Range('E18').Value = '  +2.16%  '

$ws.Range('D19').Value = '''7.78'
$ws.Range('E19').Value = '  -3.16%  '

$ws.Range('D20').Value = '''10.95'
$ws.Range('E20').Value = '  +3.22%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''327.61'
$ws.Range('E21').Value = '  +1.56%  '

$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = '''4.13'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = '''2.01'
$ws.Range('E23').Value = '  +11.57%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '''1.00'
$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('B25').Value = 'BabyDogeCoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D25').Value = '0.0₆0626'
$ws.Range('E25').Value = '  +129.43%  '

$ws.Range('D26').Value = '''65.55'
$ws.Range('E26').Value = '  +1.27%  '

$ws.Range('D27').Value = '''644.25'
$ws.Range('E27').Value = '  +14.35%  '

$ws.Range('D28').Value = '''1.17'
$ws.Range('E28').Value = '  +16.88%  '

$ws.Range('D29').Value = '''8.45'
$ws.Range('E29').Value = '  +5.84%  '

$ws.Range('D30').Value = '0.0₃0982'
$ws.Range('E30').Value = '  +5.47%  '

$ws.Range('D31').Value = '2.567.77'

$ws.Range('D32').Value = '''8.20'
$ws.Range('E32').Value = '  +2.20%  '

$ws.Range('D33').Value = '''1.43'
$ws.Range('E33').Value = '  +8.05%  '

$ws.Range('D34').Value = '''1.88'
$ws.Range('E34').Value = '  +4.28%  '

$ws.Range('D35').Value = '''0.140'
$ws.Range('E35').Value = '  +6.49%  '

$ws.Range('D36').Value = '''1.49'
$ws.Range('E36').Value = '  +2.97%  '

$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('E38').Value = '  +3.51%  '

$ws.Range('D39').Value = '''5.50'
$ws.Range('E39').Value = '  +7.32%  '

$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '''0.374'
$ws.Range('E40').Value = '  +1.26%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''152.59'
$ws.Range('E41').Value = '  +0.25%  '

$ws.Range('D42').Value = '''18.62'
$ws.Range('E42').Value = '  +2.09%  '

$ws.Range('D43').Value = '''2.71'
$ws.Range('E43').Value = '  +8.79%  '

$ws.Range('E44').Value = '  +6.02%  '

$ws.Range('D45').Value = '''42.37'
$ws.Range('E45').Value = '  +1.55%  '

$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').Value = '''15.02'
$ws.Range('E47').Value = '  +28.09%  '

$ws.Range('D48').Value = '''144.84'
$ws.Range('E48').Value = '  +2.44%  '

$ws.Range('D49').Value = '''3.61'
$ws.Range('E49').Value = '  +2.56%  '

$ws.Range('D50').Value = '''20.69'
$ws.Range('E50').Value = '  +7.87%  '

$ws.Range('E51').Value = '  +3.42%  '
